# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 03:34"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6335244
$ws.Range("C4").Value = 44507
$ws.Range("D4").Value = 3575096
$ws.Range("E4").Value = 2569090
$ws.Range("G4").Value = 1094
$ws.Range("H4").Value = 191058

# Row 6 - India
$ws.Range("D6").Value = 3034887
$ws.Range("E6").Value = 829668

# Row 23 - Alemania
$ws.Range("D23").Value = 224600
$ws.Range("E23").Value = 14815

# Row 56 - Venezuela
$ws.Range("B56").Value = 49877
$ws.Range("C56").Value = 994
$ws.Range("D56").Value = 40574
$ws.Range("E56").Value = 8901
$ws.Range("G56").Value = 4
$ws.Range("H56").Value = 402

# Row 138 - Bahamas
$ws.Range("B138").Value = 2386
$ws.Range("C138").Value = 49
$ws.Range("D138").Value = 893
$ws.Range("E138").Value = 1443

# Row 174 - San Martin (Parte Holandesa)
$ws.Range("B174").Value = 504
$ws.Range("C174").Value = 22
$ws.Range("E174").Value = 183

# Row 185 - Camboya
$ws.Range("D185").Value = 272
$ws.Range("E185").Value = 2

$wb.Save()
